$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Data edits: adding cms spread cap floor -------------------------------
# Valuation Date moves from 2019-12-31 to 2020-11-30
$ws.Range("A2").Value = 44165
# End Date Netting Set moves from 2030-04-10 to 2035-01-22
$ws.Range("B2").Value = 49331
# Collateralized flag flips from "yes" to "no"
$ws.Range("F2").Value = "no"

# --- New column for the cap/floor inputs ------------------------------------
# Widen column D slightly to fit the new content (closest reachable width)
$ws.Columns.Item(4).ColumnWidth = 12

# --- Selection moves to B3 ---------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
